$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pool")

# Remove the single "attendance" column (F) — replaced below by a monthly
# attendance rate assumption (12 columns, one per month).
$ws.Range("F1").EntireColumn.Delete()

# Add 12 new columns after "dcc" (now column AC): a_1 .. a_12, each with a
# default monthly attendance-rate value of 1.
for ($i = 1; $i -le 12; $i++) {
  $col = 29 + $i
  $headerCell = $ws.Cells.Item(1, $col)
  $valueCell = $ws.Cells.Item(2, $col)
  $headerCell.NumberFormat = "General"
  $headerCell.Value2 = "a_" + $i
  $valueCell.NumberFormat = "General"
  $valueCell.Value2 = 1
}

# The first new header cell (a_1) keeps the General format it was entered
# with; the remaining headers (a_2 .. a_12) pick up the surrounding Text
# format used by the other header cells in row 1.
$ws.Range("AE1:AO1").NumberFormat = "@"

for ($c = 1; $c -le 41; $c++) {
  $h = $ws.Cells.Item(1, $c).Value2
  $v = $ws.Cells.Item(2, $c).Value2
  Write-Host ($c.ToString() + ": " + $h + " = " + $v)
}
